$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B update
$ws.Range("B3").Value = -1182.3

# Column C updates
$ws.Range("C2").Value = -796.4
$ws.Range("C3").Value = -278.4
$ws.Range("C4").Value = -686.6
$ws.Range("C5").Value = -909.1
$ws.Range("C7").Value = -604
$ws.Range("C9").Value = -1059.8
$ws.Range("C10").Value = -647.8
$ws.Range("C11").Value = -199.3
$ws.Range("C12").Value = -434.8
$ws.Range("C13").Value = -108.7
$ws.Range("C14").Value = -355.8
$ws.Range("C15").Value = -320.2
$ws.Range("C16").Value = 718.9
$ws.Range("C17").Value = 979.3
$ws.Range("C18").Value = 1092.4
$ws.Range("C19").Value = 1686.3
$ws.Range("C20").Value = 636.7
$ws.Range("C21").Value = 292.9
$ws.Range("C22").Value = 396.3
$ws.Range("C24").Value = -152.6
